$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the income coefficient rows (A7:A14) from incomeK..incomeR
# to income2..income9, matching the updated regression output labels.
$ws.Range("A7").Value  = "income2"
$ws.Range("A8").Value  = "income3"
$ws.Range("A9").Value  = "income4"
$ws.Range("A10").Value = "income5"
$ws.Range("A11").Value = "income6"
$ws.Range("A12").Value = "income7"
$ws.Range("A13").Value = "income8"
$ws.Range("A14").Value = "income9"

# Select the edited range, mirroring the workbook's saved selection state.
$ws.Range("A7:A14").Select()
